# remplacement - par /
# The old schema-correspondance cells used " - " as a separator between
# the "old" field/table names that map to several "new" ones. The commit
# replaces that separator with " / " in every cell that uses it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targets = @("B10", "B12", "B17", "A28", "C33", "D33")

foreach ($addr in $targets) {
    $cell = $ws.Range($addr)
    $orig = $cell.Value()
    $cell.Value = $orig -replace " - ", " / "
}

# The active selection also moved one column to the left (F20 -> E20).
$ws.Range("E20").Select()
